$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in D1: remove the leading space from "Net Income"
$ws.Range("D1").Value = "Net Income"

# Update the selected cell/range to D1 (as saved in the sheet view)
$ws.Range("D1").Select()
